# Apply updated cryptocurrency price/volume data (and two coin-rank swaps)
# to match the refreshed "cryptos" export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.866.61"
$ws.Range("E2").Value = "  -0.03%  "

# Row 3
$ws.Range("D3").Value = "3.799.45"
$ws.Range("E3").Value = "  -1.51%  "

# Row 4
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.19%  "

# Row 5
$ws.Range("D5").Value = "'598.43"
$ws.Range("E5").Value = "  -0.01%  "

# Row 6
$ws.Range("D6").Value = "'169.27"
$ws.Range("E6").Value = "  +1.36%  "

# Row 7
$ws.Range("D7").Value = "3.798.92"
$ws.Range("E7").Value = "  -1.49%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").Value = "'0.525"
$ws.Range("E9").Value = "  -0.27%  "

# Row 10
$ws.Range("E10").Value = "  +0.68%  "

# Row 11
$ws.Range("D11").Value = "'6.49"
$ws.Range("E11").Value = "  +1.82%  "

# Row 12
$ws.Range("D12").Value = "'0.458"
$ws.Range("E12").Value = "  +0.54%  "

# Row 13
$ws.Range("D13").Value = "'0.0000280"
$ws.Range("E13").Value = "  +13.18%  "

# Row 14
$ws.Range("D14").Value = "'36.67"
$ws.Range("E14").Value = "  -0.46%  "

# Row 15
$ws.Range("D15").Value = "4.432.75"
$ws.Range("E15").Value = "  -1.14%  "

# Row 16
$ws.Range("D16").Value = "3.793.99"
$ws.Range("E16").Value = "  -1.72%  "

# Row 17
$ws.Range("D17").Value = "67.891.66"
$ws.Range("E17").Value = "  -0.11%  "

# Row 18
$ws.Range("D18").Value = "'18.16"
$ws.Range("E18").Value = "  +0.23%  "

# Row 19
$ws.Range("D19").Value = "'7.37"
$ws.Range("E19").Value = "  +0.20%  "

# Row 20
$ws.Range("E20").Value = "  +0.44%  "

# Row 21
$ws.Range("D21").Value = "'10.76"
$ws.Range("E21").Value = "  -1.63%  "

# Row 22
$ws.Range("D22").Value = "'467.48"
$ws.Range("E22").Value = "  +0.59%  "

# Row 23
$ws.Range("D23").Value = "'0.723"
$ws.Range("E23").Value = "  -0.22%  "

# Row 24
$ws.Range("D24").Value = "'0.0000154"
$ws.Range("E24").Value = "  -5.90%  "

# Row 25
$ws.Range("D25").Value = "'83.08"
$ws.Range("E25").Value = "  -0.03%  "

# Row 26
$ws.Range("D26").Value = "'2.25"
$ws.Range("E26").Value = "  -0.19%  "

# Row 27
$ws.Range("D27").Value = "'12.05"
$ws.Range("E27").Value = "  -0.26%  "

# Row 28
$ws.Range("E28").Value = "  +1.63%  "

# Row 30
$ws.Range("D30").Value = "'2.92"
$ws.Range("E30").Value = "  -1.10%  "

# Row 31
$ws.Range("D31").Value = "3.945.74"
$ws.Range("E31").Value = "  -1.57%  "

# Row 32
$ws.Range("D32").Value = "'7.69"
$ws.Range("E32").Value = "  -0.30%  "

# Row 33
$ws.Range("D33").Value = "'2.28"
$ws.Range("E33").Value = "  -1.50%  "

# Row 34
$ws.Range("D34").Value = "'30.76"
$ws.Range("E34").Value = "  -1.05%  "

# Row 35
$ws.Range("D35").Value = "'9.32"
$ws.Range("E35").Value = "  -0.22%  "

# Row 36
$ws.Range("D36").Value = "3.762.15"
$ws.Range("E36").Value = "  -1.85%  "

# Row 37
$ws.Range("E37").Value = "  +15.51%  "

# Row 38
$ws.Range("D38").Value = "'0.106"
$ws.Range("E38").Value = "  +2.34%  "

# Row 39
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.139"
$ws.Range("E39").Value = "  -0.26%  "

# Row 40
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").Value = "'1.01"
$ws.Range("E40").Value = "  -1.07%  "

# Row 41
$ws.Range("D41").Value = "'5.92"
$ws.Range("E41").Value = "  +0.40%  "

# Row 42
$ws.Range("D42").Value = "'0.997"
$ws.Range("E42").Value = "  -0.27%  "

# Row 43
$ws.Range("D43").Value = "'0.315"
$ws.Range("E43").Value = "  +1.11%  "

# Row 44
$ws.Range("E44").Value = "  -0.02%  "

# Row 45
$ws.Range("D45").Value = "'8.76"
$ws.Range("E45").Value = "  +2.94%  "

# Row 46
$ws.Range("D46").Value = "'1.97"
$ws.Range("E46").Value = "  -0.15%  "

# Row 47
$ws.Range("D47").Value = "'0.000291"
$ws.Range("E47").Value = "  +5.49%  "

# Row 48
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "'407.27"
$ws.Range("E48").Value = "  -4.86%  "

# Row 49
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'46.38"
$ws.Range("E49").Value = "  -1.84%  "

# Row 50
$ws.Range("D50").Value = "'141.20"
$ws.Range("E50").Value = "  -1.70%  "

# Row 51
$ws.Range("D51").Value = "'0.0356"
$ws.Range("E51").Value = "  +0.66%  "
